$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trafo")
$ws.Activate()

$data = @(
    @(1, 5, 8, 0, 0.0001, 0, 1, 0, 1),
    @(17, 30, 0, 0.0001, 0, 1, 0, 1),
    @(25, 26, 0, 0.0001, 0, 1, 0, 1),
    @(37, 38, 0, 0.0001, 0, 1, 0, 1),
    @(59, 63, 0, 0.0001, 0, 1, 0, 1),
    @(61, 64, 0, 0.0001, 0, 1, 0, 1),
    @(65, 66, 0, 0.0001, 0, 1, 0, 1),
    @(68, 69, 0, 0.0001, 0, 1, 0, 1),
    @(80, 81, 0, 0.0001, 0, 1, 0, 1)
)

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0.0001
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 1

$ws.Range("A3").Formula = "=A2+1"
$ws.Range("B3").Value = 17
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0.0001
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 1

$ws.Range("A4").Formula = "=A3+1"
$ws.Range("B4").Value = 25
$ws.Range("C4").Value = 26
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.0001
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 1

$ws.Range("A5").Formula = "=A4+1"
$ws.Range("B5").Value = 37
$ws.Range("C5").Value = 38
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.0001
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 1

$ws.Range("A6").Formula = "=A5+1"
$ws.Range("B6").Value = 59
$ws.Range("C6").Value = 63
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.0001
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 1

$ws.Range("A7").Formula = "=A6+1"
$ws.Range("B7").Value = 61
$ws.Range("C7").Value = 64
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.0001
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 1

$ws.Range("A8").Formula = "=A7+1"
$ws.Range("B8").Value = 65
$ws.Range("C8").Value = 66
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.0001
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 1

$ws.Range("A9").Formula = "=A8+1"
$ws.Range("B9").Value = 68
$ws.Range("C9").Value = 69
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.0001
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 1

$ws.Range("A10").Formula = "=A9+1"
$ws.Range("B10").Value = 80
$ws.Range("C10").Value = 81
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.0001
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 1

$ws.Range("A2:I10").Select()
